# "Generate Report for Handback" - mark the a.md localization as handed back
# for both zh-cn and de-de targets, and flip the Overview status text from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$linkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbdbf6949f58178faa2efdb33df005f3baa76ea3/e2e/a.md"

# ---------------------------------------------------------------------------
# Overview sheet: per-language status cells (E/F for rows 2 and 3)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------------
# Helper: stamp the "handback" columns (Latest Target File / Latest Handback
# File / Latest Handback DateTime = I/J/K) for a language sheet, and also
# flip the row-level Status column (C) to match the new status text.
# ---------------------------------------------------------------------------
function Set-HandbackColumns($ws, $xlfName, $handbackTime) {
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    $ws.Range("I2").Value = "a.md"
    $ws.Range("J2").Value = $xlfName
    $ws.Range("K2").Value = $handbackTime

    $ws.Range("I3").Value = "a.md"
    $ws.Range("J3").Value = $xlfName
    $ws.Range("K3").Value = $handbackTime

    # Register the actual hyperlinks for the newly-populated target files.
    $ws.Hyperlinks.Add($ws.Range("I2"), $linkUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $linkUrl, "", "", "a.md")

    # Match the look of the existing hyperlink cells (single underline,
    # cornflowerblue) on the two new "Latest Target File" links.
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackColumns $zhcn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-27 18:44:47"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-HandbackColumns $dede "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-27 18:44:53"

# ---------------------------------------------------------------------------
# Widen columns that now hold the longer status / handback text so the
# report is readable (mirrors Excel's column auto-fit after the edit).
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(10).ColumnWidth = 40
